# passage m1 a m2
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new columns before the current "Parcours" column (C),
# so the header order becomes: Nom, Prenom, Promo, TypeAnnee, Parcours, ...
$ws.Range("C1:D1").EntireColumn.Insert()

$ws.Range("C1").Value = "Promo"
$ws.Range("D1").Value = "TypeAnnée"

# New student rows (post M1->M2 transition)
$data = @(
    @("ABBE", "TRISTAN", 2025, "M2", "GPhy", "", "apprentissage", "SANOFI", "Gentilly (94)", "BELLOCQ", "", "GENIET"),
    @("KONE", "YACOUBA", 2025, "M1", "GPhy", "", "pro", "LABORATOIRE XLIM", "Poitiers", "BOURDON", "", "URRUTY"),
    @("MONTBULEAU--GENTELET", "TITOUAN", 2025, "M1", "GPhy", "", "apprentissage", "MAAT PHARMA", "Lyon", "BERGÉ", "0617421317", "URRUTY"),
    @("NIGGEL", "THIBAULT", 2025, "M2", "GPhy", "", "stage", "EVALU CONSEIL", "Paris / à distance", "MONNé", "0615097890", "URRUTY"),
    @("SAVANY", "THIBAULT", 2025, "M2", "GPhy", "", "pro", "", "", "", "", "")
)

# Column K holds phone numbers that must stay text (leading zeros)
$ws.Range("K4:K5").NumberFormat = "@"

$rowIndex = 2
foreach ($row in $data) {
    $colIndex = 1
    foreach ($val in $row) {
        if ($val -ne "") {
            $ws.Cells.Item($rowIndex, $colIndex).Value = $val
        }
        $colIndex++
    }
    $rowIndex++
}
